# Parallel tree execution output: add the memory-usage samples that were
# previously missing for rows 3 and 4 of columns I (tree depth) and J
# (tree breadth) on the "Memory Usage" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Memory Usage")

$ws.Range("I3").Value = 3120096.0
$ws.Range("J3").Value = 2413312.0

$ws.Range("I4").Value = 7197608.0
$ws.Range("J4").Value = 1313128.0
